# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.326.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6279"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -2.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2894"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.840.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.987"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6780"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.37%  "

# Row 15
$ws.Range("E15").Value = "  -2.74%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.258"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.328.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.428"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.473"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1353"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06645"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.88%  "

# Row 29
$ws.Range("E29").Value = "  +3.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.071"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.36%  "

# Row 32
$ws.Range("E32").Value = "  +0.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.836"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.37%  "

# Row 34
$ws.Range("E34").Value = "  -1.50%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6933"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.580"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01859"
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.823"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.52%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.246.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9378"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.51%  "

# Row 42
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.005.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.035"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "

# Row 47
$ws.Range("E47").Value = "  +2.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.009"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1150"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3899"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.51%  "
